# Estado de Cuenta - actualizacion de base de datos de trabajadores en mora
# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Move the signature/footer block (currently rows 23-24, each split into
#    a B:C merged cell and an H:J merged cell) down to rows 30-31, so the
#    new, larger worker table (which now needs through row 25) has room.
#    Copy each merged block individually (copying the whole B:J swath at
#    once pulls in the plain column styles for the untouched D:G cells).
# ---------------------------------------------------------------------------
$ws.Range("B23:C23").Copy($ws.Range("B30:C30"))
$ws.Range("H23:J23").Copy($ws.Range("H30:J30"))
$ws.Range("B24:C24").Copy($ws.Range("B31:C31"))
$ws.Range("H24:J24").Copy($ws.Range("H31:J31"))
$excel.CutCopyMode = $false

$ws.Range("B23:C23").UnMerge()
$ws.Range("H23:J23").UnMerge()
$ws.Range("B24:C24").UnMerge()
$ws.Range("H24:J24").UnMerge()

# Drop the now-duplicated text that used to live in rows 23-24 (those rows
# become worker/mora data rows below).
$ws.Range("B23:J24").ClearContents()

# ---------------------------------------------------------------------------
# 2) Replicate the worker-row formatting down through the new rows.
#    Rows 16/17 carry the "interior" row style; row 18 (before this edit)
#    carried the "last row" (heavier bottom border) style.
#    After the edit: rows 16-24 use the interior style, row 25 uses the
#    "last row" style.
# ---------------------------------------------------------------------------
$ws.Range("B18:J18").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)  # xlPasteFormats  (grab the "last row" style before it's overwritten)
$excel.CutCopyMode = $false

$ws.Range("B17:J17").Copy()
$ws.Range("B18:J24").PasteSpecial(-4122)  # xlPasteFormats (interior style for all the new middle rows)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Write the new worker / mora table (rows 16-25).
#    4 distinct workers, 7 total periods of the recurring worker.
# ---------------------------------------------------------------------------
$rows = @(
  @{ R=16; Tipo="CC"; Doc="73147370";   Nombre="LAUREANO GARRIDO MIRANDA";        Periodo="2507"; Mora=160000; Salario=4000000 },
  @{ R=17; Tipo="CC"; Doc="1043962336"; Nombre="ISAURA ANDREA MARRUGO SUAREZ";    Periodo="2507"; Mora=56940;  Salario=1423500 },
  @{ R=18; Tipo="CC"; Doc="1001970934"; Nombre="MARIA CAMILA VILLAR LOPEZ";       Periodo="2507"; Mora=52000;  Salario=1300000 },
  @{ R=19; Tipo="CC"; Doc="1001970934"; Nombre="MARIA CAMILA VILLAR LOPEZ";       Periodo="2506"; Mora=52000;  Salario=1300000 },
  @{ R=20; Tipo="CC"; Doc="1001970934"; Nombre="MARIA CAMILA VILLAR LOPEZ";       Periodo="2505"; Mora=52000;  Salario=1300000 },
  @{ R=21; Tipo="CC"; Doc="1001970934"; Nombre="MARIA CAMILA VILLAR LOPEZ";       Periodo="2504"; Mora=52000;  Salario=1300000 },
  @{ R=22; Tipo="CC"; Doc="1001970934"; Nombre="MARIA CAMILA VILLAR LOPEZ";       Periodo="2503"; Mora=52000;  Salario=1300000 },
  @{ R=23; Tipo="CC"; Doc="1001970934"; Nombre="MARIA CAMILA VILLAR LOPEZ";       Periodo="2502"; Mora=52000;  Salario=1300000 },
  @{ R=24; Tipo="CC"; Doc="1001970934"; Nombre="MARIA CAMILA VILLAR LOPEZ";       Periodo="2501"; Mora=52000;  Salario=1300000 },
  @{ R=25; Tipo="CC"; Doc="1047506058"; Nombre="MARIA ALEJANDRA BELTRAN GOMEZ";   Periodo="2507"; Mora=56940;  Salario=1423500 }
)

foreach ($row in $rows) {
  $r = $row.R
  $ws.Range("B$r").Value2 = $row.Tipo
  $ws.Range("C$r").Value2 = $row.Doc
  $ws.Range("D$r").Value2 = $row.Nombre
  $ws.Range("E$r").Value2 = $row.Periodo
  $ws.Range("F$r").Value2 = $row.Mora
  $ws.Range("G$r").Value2 = $row.Salario
}

# ---------------------------------------------------------------------------
# 4) Re-merge the footer blocks at their new location.
# ---------------------------------------------------------------------------
$ws.Range("B30:C30").Merge()
$ws.Range("H30:J30").Merge()
$ws.Range("B31:C31").Merge()
$ws.Range("H31:J31").Merge()

# ---------------------------------------------------------------------------
# 5) Update the summary header cells: total mora, worker count, period count.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 637880
$ws.Range("C13").Value2 = 4
$ws.Range("F13").Value2 = 7

Write-Output "Estado de cuenta actualizado: 4 trabajadores, 7 periodos, mora total 637880"
